$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.210.39"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "2.642.06"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.32%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -1.56%  "

$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("D14").Value = "3.123.97"
$ws.Range("E14").Value = "  +0.08%  "

$ws.Range("E15").Value = "  -2.86%  "

$ws.Range("D16").Value = "68.174.63"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "2.686.63"
$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "358.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("E22").Value = "  -3.62%  "

$ws.Range("E23").Value = "  -0.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("D27").Value = "2.815.48"
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "556.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.32%  "

$ws.Range("E31").Value = "  -2.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("E40").Value = "  -2.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.73%  "

$ws.Range("E43").Value = "  -7.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "156.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("E48").Value = "  -3.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0773"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.567"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
